$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F (dSF)
$updates = @{
    3  = 1
    4  = 5
    5  = 4
    6  = 3
    7  = 9
    8  = 0
    9  = -8
    10 = -3
    11 = 3
    12 = 3
    13 = -1
    14 = -3
    15 = -2
    16 = 0
    17 = 11
    18 = 5
    19 = 5
    20 = -2
    21 = -5
    23 = -1
    24 = -3
    25 = -8
    26 = -4
    27 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
